# Update scraped_at timestamps (column K) on the "snapshot" worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("snapshot")

$timestamps = @{
    2 = "2025-12-18T04:37:57.313773+00:00"
    3 = "2025-12-18T04:37:57.313793+00:00"
    4 = "2025-12-18T04:37:59.447955+00:00"
    5 = "2025-12-18T04:37:59.447990+00:00"
    6 = "2025-12-18T04:37:59.448011+00:00"
    7 = "2025-12-18T04:38:01.546874+00:00"
    8 = "2025-12-18T04:38:03.611056+00:00"
    9 = "2025-12-18T04:38:05.327264+00:00"
    10 = "2025-12-18T04:38:05.327287+00:00"
    11 = "2025-12-18T04:38:07.427180+00:00"
    12 = "2025-12-18T04:38:11.175469+00:00"
    13 = "2025-12-18T04:38:11.175492+00:00"
    14 = "2025-12-18T04:38:13.331824+00:00"
    15 = "2025-12-18T04:38:15.430629+00:00"
    16 = "2025-12-18T04:38:17.532109+00:00"
    17 = "2025-12-18T04:38:19.178451+00:00"
    18 = "2025-12-18T04:38:19.178479+00:00"
    19 = "2025-12-18T04:38:19.178495+00:00"
    20 = "2025-12-18T04:38:19.178511+00:00"
    21 = "2025-12-18T04:38:20.834729+00:00"
    22 = "2025-12-18T04:38:20.834761+00:00"
    23 = "2025-12-18T04:38:22.507016+00:00"
    24 = "2025-12-18T04:38:22.507034+00:00"
    25 = "2025-12-18T04:38:22.507042+00:00"
    26 = "2025-12-18T04:38:24.701313+00:00"
    27 = "2025-12-18T04:38:24.701343+00:00"
    28 = "2025-12-18T04:38:26.944959+00:00"
    29 = "2025-12-18T04:38:26.944988+00:00"
    30 = "2025-12-18T04:38:26.945005+00:00"
    31 = "2025-12-18T04:38:29.084227+00:00"
    32 = "2025-12-18T04:38:31.297003+00:00"
    33 = "2025-12-18T04:38:31.297029+00:00"
    34 = "2025-12-18T04:38:35.204793+00:00"
    35 = "2025-12-18T04:38:35.204814+00:00"
    36 = "2025-12-18T04:38:37.266625+00:00"
    37 = "2025-12-18T04:38:37.266653+00:00"
}

foreach ($row in $timestamps.Keys) {
    $ws.Range("K$row").Value = $timestamps[$row]
}

